$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.408222913742065
$ws.Range("B1").Value = 1.981092810630798
$ws.Range("C1").Value = 2.941610336303711
$ws.Range("D1").Value = 4.80394983291626
$ws.Range("E1").Value = 0.937427282333374
